$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-day column CM, dated 2025-12-09 (serial 46000), appended
# right after the current last column CL (2025-12-05 / serial 45996).

# 1) Header date cell (row 1)
$ws.Range("CM1").Value = 46000
$ws.Range("CL1").Copy()
$ws.Range("CM1").PasteSpecial(-4122)   # xlPasteFormats - reuse CL's date style

# 2) Attendance mark for each player row. Row 12 gets no new cell (that
#    player's data already stops before column CL/CM) and row 21 gets a
#    styled-but-empty cell (that player's series already ended earlier).
$values = @{
    2  = "P";
    3  = "R";
    4  = "P";
    5  = "P";
    6  = "P";
    7  = "P";
    8  = "P";
    9  = "P";
    10 = "P";
    11 = "P";
    13 = "B";
    14 = "P";
    15 = "P";
    16 = "P";
    17 = "P";
    18 = "P";
    19 = "P";
    20 = "P";
    22 = "P";
    23 = "B";
    24 = "P";
    25 = "RH";
    26 = "P";
    27 = "P";
    28 = "P";
    29 = "P"
}

foreach ($r in $values.Keys) {
    $ws.Range("CM$r").Value = $values[$r]
}

foreach ($r in 1..29) {
    if ($r -eq 1 -or $r -eq 12) { continue }
    $ws.Range("CL$r").Copy()
    $ws.Range("CM$r").PasteSpecial(-4122)   # xlPasteFormats
}

# 3) Matches where the user last clicked after entering the new column.
[void]$ws.Range("CO27").Select()
